# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (column E, rows 16-25) currently lists periods
# in descending order (1710 .. 1701). This update refreshes the account
# statement data so the periods run in ascending order (1701 .. 1710),
# keeping every other value/format on the row untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1701", "1702", "1703", "1704", "1705", "1706", "1707", "1708", "1709", "1710")

$firstRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}
